$d = $word.ActiveDocument

# The document has Pearson logo images (PNG) embedded in the footers and
# BTec logo images (JPG) embedded in the headers. This commit swaps the
# display/"name" metadata recorded on each picture:
#   footers: PearsonLogo pictures  image1.png -> image2.png
#   headers: BTec_Logo-Orange pictures  image2.jpg -> image1.jpg
# The pictures themselves (their binary data / relationship) are untouched;
# only the shape's Name metadata changes.

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
                $shp = $h.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
                $shp = $f.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
